$d = $word.ActiveDocument

# Build the combined plain text for all new paragraphs, each ended by a paragraph mark.
# Using InsertBefore first establishes correct paragraph properties (inherited from the
# paragraph that used to be first), matching the existing "plain" note paragraphs exactly.
$combinedText = "16/09/25 Presentación 08`rLos gráficos deben tener título tanto para el gráfico en sí como para los ejes.`rHay riesgos que no tienen plan de mitigación (presentamos 3, quedaban 5)`rNo pusimos los Casos de Uso`rPara cualquier tipo de plan, no asignar a las tareas plazos que sean más grandes que la iteración. Asignar momento específico de la iteración en el que se desarrollará la tarea.`rLas métricas deben tener un valor planificado y un valor ejecutado para tener margen de comparación. Además son por iteración.`rLa estimación está mal realizada. El valor que brinda la metodología que utilizamos, es del 40%. El 60% restante es el que debemos agregar nosotros, o algo así entendí. `rNo pusimos el gráfico de estimación que se pide. Diapositivas obligatorias como riesgos, estimación y CU, deben ocupar una sola diapositiva y deben ser un resumen rápido no extremadamente detallado.`rEl gráfico de estimación debe tener la fecha estimada de finalización puesta.`rLa estimación se io sobreestimada, aparentemente, porque los casos de uso están sobreestimados ahre no se como decirlo, pero puede que hayan más casos de uso de los que realmente deberían. Con respecto a esto después nos dijeron que los crud son un solo caso de uso.`r"
$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertBefore($combinedText)

# Word always serialises <w:t> runs it authors with xml:space="preserve" (matching the rest
# of this note-taking document). The plain InsertBefore call above only adds that attribute
# when the text has leading/trailing whitespace, so patch each newly created run back in via
# InsertXML (scoped to the text only, so paragraph formatting/<w:pPr> stays untouched).
$newParagraphCount = 10
$texts = @(
  "16/09/25 Presentación 08",
  "Los gráficos deben tener título tanto para el gráfico en sí como para los ejes.",
  "Hay riesgos que no tienen plan de mitigación (presentamos 3, quedaban 5)",
  "No pusimos los Casos de Uso",
  "Para cualquier tipo de plan, no asignar a las tareas plazos que sean más grandes que la iteración. Asignar momento específico de la iteración en el que se desarrollará la tarea.",
  "Las métricas deben tener un valor planificado y un valor ejecutado para tener margen de comparación. Además son por iteración.",
  "La estimación está mal realizada. El valor que brinda la metodología que utilizamos, es del 40%. El 60% restante es el que debemos agregar nosotros, o algo así entendí. ",
  "No pusimos el gráfico de estimación que se pide. Diapositivas obligatorias como riesgos, estimación y CU, deben ocupar una sola diapositiva y deben ser un resumen rápido no extremadamente detallado.",
  "El gráfico de estimación debe tener la fecha estimada de finalización puesta.",
  "La estimación se io sobreestimada, aparentemente, porque los casos de uso están sobreestimados ahre no se como decirlo, pero puede que hayan más casos de uso de los que realmente deberían. Con respecto a esto después nos dijeron que los crud son un solo caso de uso."
)

for ($i = 0; $i -lt $newParagraphCount; $i++) {
  $p = $d.Paragraphs.Item($i + 1)
  $pStart = $p.Range.Start
  $pEnd = $p.Range.End - 1
  $textRange = $d.Range($pStart, $pEnd)
  $escaped = $texts[$i] -replace "&", "&amp;" -replace "<", "&lt;" -replace ">", "&gt;"
  $runXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes" ?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  $textRange.InsertXML($runXml)
}
